# aproda-testreport.xlsx update
#
# Summary of intent (from commit message + diff):
#  - "import budget name from different column": the "Aufwände gesamt" sheet's
#    column F (previously blank/space placeholder under the "Subgruppe" header)
#    now carries the literal label "Budget" for every data row.
#  - Print titles on "Aufwände gesamt" extended from rows 1:3 to 1:4 (matching
#    the other two sheets, which already repeat rows 1-4 on each printed page).
#  - Cosmetic: sheet tab colors normalized to opaque white (alpha byte fixed
#    from 00 to FF) on all three sheets.
#  - The active selection follows the relocated "Budget" cell: F11 instead of
#    the old E11 on the "Aufwände gesamt" sheet (the sheet that was active/
#    selected when the workbook was last saved).

$wb = $excel.ActiveWorkbook

$wsArbeitspakete   = $wb.Worksheets.Item(1)   # "Arbeitspakete"
$wsArbeitspaketePM = $wb.Worksheets.Item(2)   # "Arbeitspakete PM"
$wsAufwaende       = $wb.Worksheets.Item(3)   # "Aufwände gesamt"

# --- content edit: budget name moved into column F (rows 4-11) -------------
$wsAufwaende.Range("F4:F11").Value = "Budget"

# --- print titles: repeat rows 1-4 (was 1-3) on "Aufwaende gesamt" ---------
$wsAufwaende.PageSetup.PrintTitleRows = '$1:$4'

# --- tab colors: opaque white on every sheet --------------------------------
$wsArbeitspakete.Tab.Color   = 16777215
$wsArbeitspaketePM.Tab.Color = 16777215
$wsAufwaende.Tab.Color       = 16777215

# --- selection follows the relocated "Budget" column (F11) -----------------
$wsAufwaende.Activate()
$wsAufwaende.Range("F11").Select()
